# Update LCOE comparison workbook:
#  - add two new operating-point runs (15op and 20op, with/without PV) to the
#    "Comparison_operating_points" sheet/table
#  - make "Comparison_operating_points" the active sheet/tab (selecting F31)
#  - drop the stray highlight style on the "PV" helper-column formula cells of
#    the "Comparison_sensitivities" sheet

$wb = $excel.ActiveWorkbook

$wsOperating = $wb.Worksheets.Item("Comparison_operating_points")
$wsSensitivities = $wb.Worksheets.Item("Comparison_sensitivities")

# --- 1. Extend Table1 on the operating-points sheet and fill in the new rows ---

$table = $wsOperating.ListObjects.Item("Table1")
$table.Resize($wsOperating.Range("A1:I25"))

$newRuns = @(
    @{ Row = 22; Name = "base_15op";    LcoeMWh = 249.73614149918711; LcoeT = 1380.485893287173;  Invest = 288736929.04657698; Annual = 14767064.360440319; EnergyMWh = 176888.9279999999; EnergyT = 32000.00707537687;  Pcf = 9.8181474074492936 },
    @{ Row = 23; Name = "base_15op_PV"; LcoeMWh = 196.44983468472191; LcoeT = 1085.931030618324;  Invest = 288736929.04657698; Annual = 5341306.6709504724;  EnergyMWh = 176888.9279999999; EnergyT = 32000.00707537687;  Pcf = 9.8181474074492936 },
    @{ Row = 24; Name = "base_20op";    LcoeMWh = 249.59496612046431; LcoeT = 1379.7055071658999; Invest = 288563658.52300167; Annual = 14759739.98456678;  EnergyMWh = 176888.9279999143; EnergyT = 32000.007075361391; Pcf = 9.8181474074492936 },
    @{ Row = 25; Name = "base_20op_PV"; LcoeMWh = 196.2722336916267;  LcoeT = 1084.949291795381;  Invest = 288563658.52300167; Annual = 5327539.0072034188;  EnergyMWh = 176888.9279999143; EnergyT = 32000.007075361391; Pcf = 9.8181474074492936 }
)

foreach ($run in $newRuns) {
    $r = $run.Row
    $wsOperating.Range("A$r").Value = $run.Name
    $wsOperating.Range("B$r").Formula = '=IF(ISNUMBER(SEARCH("PV", Table1[[#This Row],[run_name]])),"PV revenue","no PV revenue")'
    $wsOperating.Range("C$r").Value = $run.LcoeMWh
    $wsOperating.Range("D$r").Value = $run.LcoeT
    $wsOperating.Range("E$r").Value = $run.Invest
    $wsOperating.Range("F$r").Value = $run.Annual
    $wsOperating.Range("G$r").Value = $run.EnergyMWh
    $wsOperating.Range("H$r").Value = $run.EnergyT
    $wsOperating.Range("I$r").Value = $run.Pcf
}

# --- 2. Remove the fill/border style from the "PV" helper column on the
#        sensitivities sheet (cells B44 and B46:B62) ---

$wsSensitivities.Range("B44").Style = "Normal"
$wsSensitivities.Range("B46:B62").Style = "Normal"

# --- 3. Make the operating-points sheet the active tab / selection ---

$wsOperating.Activate()
$wsOperating.Range("F31").Select()
